# RegisterCommandExample.pptx edit
# Summary of change (per commit message "diagrams: added additional supporting
# images and resolve discrepancies with existing diagrams"):
#   - The single existing slide ("before/after register" diagram) is duplicated
#     so the deck now has two slides; the duplicate becomes the new slide 2
#     and keeps the original banner picture + layout untouched.
#   - The original slide (now slide 1) is reworked into a side-by-side
#     "before execution" figure: the two UI screenshots and the red
#     highlight box are shifted right to make room, the big banner picture
#     is removed, and two new italic caption text boxes are added
#     underneath ("Fig 1a. ..." / "Fig 1b. ...").

$p = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)

# EMU -> point conversion helper (1 pt = 12700 EMU); PowerPoint's COM object
# model works in points for Left/Top/Width/Height.
function EMU($v) { return $v / 12700.0 }

# --- 1. Duplicate the only slide -------------------------------------------------
# This creates slide 2 as an exact copy (same shapes/ids/positions/rels) and
# keeps the original slide in position 1, matching the new <p:sldId> order
# (original id=256 stays first, new id=257 slide follows).
$s2 = $s1.Duplicate()

# --- 2. Rework slide 1 -------------------------------------------------------

# 2a. Shift the "after register login" screenshot to the right.
$picAfter = $s1.Shapes.Item("Picture 21")
$picAfter.Left = EMU(6648547)

# 2b. Shift the red callout rounded-rectangle to the right (stays aligned to
#     the screenshot it highlights).
$rectRed = $s1.Shapes.Item("Rounded Rectangle 16")
$rectRed.Left = EMU(7736937)

# 2c. Remove the big top banner picture -- it is not part of slide 1 anymore
#     (it now only lives on the duplicated slide 2).
$banner = $s1.Shapes.Item("Picture 1")
$banner.Delete()

# 2d. Burn through a few throwaway shape-id allocations so the two new
#     caption boxes line up with the ids PowerPoint assigned them
#     (id 7 and 8) when this edit was authored.
for ($i = 1; $i -le 4; $i++) {
    $tmp = $s1.Shapes.AddTextbox(1, 0, 0, 10, 10)
    $tmp.Delete()
}

# 2e. Caption under the left-hand (before) screenshot.
$cap1Left = EMU(1027211)
$cap1Top = EMU(5538765)
$cap1Width = EMU(4536242)
$cap1Height = EMU(369332)
$cap1 = $s1.Shapes.AddShape(1, $cap1Left, $cap1Top, $cap1Width, $cap1Height)
$cap1.Name = "Rectangle 6"
$cap1.TextFrame.WordWrap = 0
$cap1.TextFrame.AutoSize = 1
$tr1 = $cap1.TextFrame.TextRange
$tr1.Text = "Fig 1a. Before execution of Register command."
$tr1.Font.Italic = $true
$tr1.Font.Name = "Times New Roman"

# 2f. Caption under the right-hand (after/expected) screenshot -- typed as
#     two runs, same as the authored deck ("Fig 1b" then the remainder).
$cap2Left = EMU(6349028)
$cap2Top = EMU(5538765)
$cap2Width = EMU(4685898)
$cap2Height = EMU(369332)
$cap2 = $s1.Shapes.AddShape(1, $cap2Left, $cap2Top, $cap2Width, $cap2Height)
$cap2.Name = "Rectangle 7"
$cap2.TextFrame.WordWrap = 0
$cap2.TextFrame.AutoSize = 1
$tr2 = $cap2.TextFrame.TextRange
$tr2.Text = "Fig 1b"
[void]$tr2.InsertAfter(". Expected UI Output(Register Command)")
$capFull2 = $cap2.TextFrame.TextRange
$capFull2.Font.Italic = $true
$capFull2.Font.Name = "Times New Roman"
